# Form the consolidated report: fill in the "Absent" column (H) values
# for the attendance rows that were previously blank/incorrect.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H12").Value = 0
